# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values on the zh-cn and de-de report sheets to reflect the re-generated
# handback report timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D5").Value = "2016-02-06 04:29:42"
$zhcn.Range("G5").Value = "2016-02-06 04:30:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D5").Value = "2016-02-06 04:29:53"
$dede.Range("G5").Value = "2016-02-06 04:30:49"
